$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Case")

# Insert a new column before column A, shifting all existing data
# (columns A:K -> B:L) just like the authored edit did.
$ws.Columns.Item(1).Insert()

# New header cell for the inserted column.
$ws.Range("A1").Value = "Automate"

# Mark every test-case block ("Y") in the new column.
$markedRows = @(3, 7, 11, 15, 19, 23)
foreach ($r in $markedRows) {
    $cell = $ws.Range("A$r")
    $cell.WrapText = $true
    $cell.Value = "Y"
}

Write-Host "done"
